$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1136922373902997
$ws.Range("C2").Value = 0.5976608816434473
$ws.Range("D2").Value = 0.5710810818331241
$ws.Range("E2").Value = 0.7556990683024057
$ws.Range("F2").Value = 0.7687573175481709
$ws.Range("G2").Value = 18
